$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37899
$ws.Range("D2").Value = 54810913
$ws.Range("C3").Value = 91283
$ws.Range("D3").Value = 133811298
$ws.Range("C4").Value = 31250
$ws.Range("D4").Value = 46280249
$ws.Range("C5").Value = 8720
$ws.Range("D5").Value = 12960563
$ws.Range("C6").Value = 2004
$ws.Range("D6").Value = 2977971
$ws.Range("C7").Value = 154
$ws.Range("D7").Value = 226093
$ws.Range("C12").Value = 41454
$ws.Range("D12").Value = 56245776
$ws.Range("C13").Value = 9685
$ws.Range("D13").Value = 14008975
$ws.Range("C14").Value = 26023
$ws.Range("D14").Value = 38161497
$ws.Range("C15").Value = 8325
$ws.Range("D15").Value = 12355324
$ws.Range("C16").Value = 2153
$ws.Range("D16").Value = 3201665
$ws.Range("C17").Value = 419
$ws.Range("D17").Value = 617623
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 51000
$ws.Range("C20").Value = 10255
$ws.Range("D20").Value = 13571841
$ws.Range("C21").Value = 13437
$ws.Range("D21").Value = 19397994
$ws.Range("C22").Value = 31745
$ws.Range("D22").Value = 46584220
$ws.Range("C23").Value = 10244
$ws.Range("D23").Value = 15227578
$ws.Range("C24").Value = 2646
$ws.Range("D24").Value = 3934182
$ws.Range("C25").Value = 508
$ws.Range("D25").Value = 756092
$ws.Range("C27").Value = 11720
$ws.Range("D27").Value = 15654310
$ws.Range("C28").Value = 7681
$ws.Range("D28").Value = 11123160
$ws.Range("C29").Value = 22568
$ws.Range("D29").Value = 33127365
$ws.Range("C30").Value = 7837
$ws.Range("D30").Value = 11663133
$ws.Range("C31").Value = 1965
$ws.Range("D31").Value = 2931999
$ws.Range("C33").Value = 29
$ws.Range("D33").Value = 43393
$ws.Range("C34").Value = 8336
$ws.Range("D34").Value = 11009451
$ws.Range("C35").Value = 3261
$ws.Range("D35").Value = 4707847
$ws.Range("C36").Value = 7864
$ws.Range("D36").Value = 11483781
$ws.Range("C37").Value = 3185
$ws.Range("D37").Value = 4720461
$ws.Range("C41").Value = 2483
$ws.Range("D41").Value = 3356485
$ws.Range("C42").Value = 17311
$ws.Range("D42").Value = 25028278
$ws.Range("C43").Value = 51281
$ws.Range("D43").Value = 75180261
$ws.Range("C44").Value = 19064
$ws.Range("D44").Value = 28317008
$ws.Range("C45").Value = 5617
$ws.Range("D45").Value = 8364355
$ws.Range("C46").Value = 1210
$ws.Range("D46").Value = 1805545
$ws.Range("C47").Value = 62
$ws.Range("D47").Value = 91068
$ws.Range("C50").Value = 16761
$ws.Range("D50").Value = 22311337
$ws.Range("C51").Value = 2033
$ws.Range("D51").Value = 2948971
$ws.Range("C52").Value = 6929
$ws.Range("D52").Value = 10185363
$ws.Range("C53").Value = 2361
$ws.Range("D53").Value = 3526418
$ws.Range("C54").Value = 756
$ws.Range("D54").Value = 1129305
$ws.Range("C57").Value = 7034
$ws.Range("D57").Value = 9668514
$ws.Range("C58").Value = 972
$ws.Range("D58").Value = 1454339
$ws.Range("C59").Value = 2450
$ws.Range("D59").Value = 3712295
$ws.Range("C60").Value = 979
$ws.Range("D60").Value = 1500001
$ws.Range("C61").Value = 335
$ws.Range("D61").Value = 518758
$ws.Range("C62").Value = 107
$ws.Range("D62").Value = 164850
$ws.Range("C64").Value = 1421
$ws.Range("D64").Value = 2023053
$ws.Range("C65").Value = 15434
$ws.Range("D65").Value = 22291846
$ws.Range("C66").Value = 44848
$ws.Range("D66").Value = 65626012
$ws.Range("C67").Value = 15741
$ws.Range("D67").Value = 23393601
$ws.Range("C68").Value = 4578
$ws.Range("D68").Value = 6818051
$ws.Range("C69").Value = 930
$ws.Range("D69").Value = 1383168
$ws.Range("C70").Value = 78
$ws.Range("D70").Value = 114330
$ws.Range("C73").Value = 15136
$ws.Range("D73").Value = 19951209
$ws.Range("C74").Value = 51896
$ws.Range("D74").Value = 75520165
$ws.Range("C75").Value = 147056
$ws.Range("D75").Value = 216652006
$ws.Range("C76").Value = 63860
$ws.Range("D76").Value = 95160999
$ws.Range("C77").Value = 20419
$ws.Range("D77").Value = 30508831
$ws.Range("C78").Value = 4843
$ws.Range("D78").Value = 7233543
$ws.Range("C79").Value = 265
$ws.Range("D79").Value = 392670
$ws.Range("C85").Value = 51216
$ws.Range("D85").Value = 69666827
$ws.Range("C86").Value = 4627
$ws.Range("D86").Value = 6703989
$ws.Range("C87").Value = 11606
$ws.Range("D87").Value = 17051819
$ws.Range("C88").Value = 3892
$ws.Range("D88").Value = 5800583
$ws.Range("C89").Value = 1348
$ws.Range("D89").Value = 2014489
$ws.Range("C90").Value = 288
$ws.Range("D90").Value = 429512
$ws.Range("C92").Value = 5
$ws.Range("D92").Value = 7500
$ws.Range("C93").Value = 5432
$ws.Range("D93").Value = 7303702
$ws.Range("C94").Value = 1603
$ws.Range("D94").Value = 2308033
$ws.Range("C95").Value = 5181
$ws.Range("D95").Value = 7630243
$ws.Range("C96").Value = 1941
$ws.Range("D96").Value = 2891437
$ws.Range("C97").Value = 690
$ws.Range("D97").Value = 1033960
$ws.Range("C98").Value = 186
$ws.Range("D98").Value = 278113
$ws.Range("C101").Value = 3574
$ws.Range("D101").Value = 4728958
$ws.Range("C102").Value = 624
$ws.Range("D102").Value = 955164
$ws.Range("C103").Value = 362
$ws.Range("D103").Value = 553480
$ws.Range("C104").Value = 132
$ws.Range("D104").Value = 199160
$ws.Range("C106").Value = 21
$ws.Range("D106").Value = 33000
$ws.Range("C107").Value = 10807
$ws.Range("D107").Value = 15678239
$ws.Range("C108").Value = 29297
$ws.Range("D108").Value = 43042608
$ws.Range("C109").Value = 9799
$ws.Range("D109").Value = 14571412
$ws.Range("C110").Value = 2696
$ws.Range("D110").Value = 4020207
$ws.Range("C111").Value = 492
$ws.Range("D111").Value = 733046
$ws.Range("C112").Value = 51
$ws.Range("D112").Value = 76500
$ws.Range("C114").Value = 9817
$ws.Range("D114").Value = 12969283
$ws.Range("C115").Value = 30591
$ws.Range("D115").Value = 44115855
$ws.Range("C116").Value = 66337
$ws.Range("D116").Value = 97079678
$ws.Range("C117").Value = 21412
$ws.Range("D117").Value = 31822131
$ws.Range("C118").Value = 6078
$ws.Range("D118").Value = 9055021
$ws.Range("C119").Value = 1129
$ws.Range("D119").Value = 1687271
$ws.Range("C120").Value = 78
$ws.Range("D120").Value = 114420
$ws.Range("C124").Value = 25921
$ws.Range("D124").Value = 34619248
$ws.Range("C125").Value = 36143
$ws.Range("D125").Value = 52165127
$ws.Range("C126").Value = 77017
$ws.Range("D126").Value = 112621101
$ws.Range("C127").Value = 23903
$ws.Range("D127").Value = 35475356
$ws.Range("C128").Value = 6408
$ws.Range("D128").Value = 9522738
$ws.Range("C129").Value = 1242
$ws.Range("D129").Value = 1847411
$ws.Range("C130").Value = 59
$ws.Range("D130").Value = 86728
$ws.Range("C133").Value = 31895
$ws.Range("D133").Value = 42348292
$ws.Range("C134").Value = 13281
$ws.Range("D134").Value = 19223580
$ws.Range("C135").Value = 32417
$ws.Range("D135").Value = 47611217
$ws.Range("C136").Value = 11495
$ws.Range("D136").Value = 17079792
$ws.Range("C137").Value = 2964
$ws.Range("D137").Value = 4418714
$ws.Range("C141").Value = 10845
$ws.Range("D141").Value = 14460585
$ws.Range("C142").Value = 35204
$ws.Range("D142").Value = 50840354
$ws.Range("C143").Value = 81493
$ws.Range("D143").Value = 119395873
$ws.Range("C144").Value = 24404
$ws.Range("D144").Value = 36258306
$ws.Range("C145").Value = 6412
$ws.Range("D145").Value = 9567567
$ws.Range("C146").Value = 1439
$ws.Range("D146").Value = 2140730
$ws.Range("C147").Value = 81
$ws.Range("D147").Value = 121130
$ws.Range("C149").Value = 29269
$ws.Range("D149").Value = 39475444
